$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60; this shifts the existing rows 60-87 down to 61-88,
# preserving all of their data and formatting.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new record.
$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C60").Value = "Arica y Parinacota"
$ws.Range("D60").Value = 44609
$ws.Range("D60").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E60").Value = 15
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100102
$ws.Range("H60").Value = "Cítricos"
$ws.Range("I60").Value = 100102004
$ws.Range("J60").Value = "Mandarina"
$ws.Range("K60").Value = "Murcott"
$ws.Range("L60").Value = "Tercera"
$ws.Range("M60").Value = 250
$ws.Range("N60").Value = 13000
$ws.Range("O60").Value = 14000
$ws.Range("P60").Value = 13500
$ws.Range("Q60").Value = "`$/caja 20 kilos"
$ws.Range("R60").Value = "Región de Coquimbo"
$ws.Range("S60").Value = 675
$ws.Range("T60").Value = 20
